$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "D1_USD" (sheet1.xml) — append a new day (45264) as row 76, and
# backfill C77 with the prediction that arrived for the following day.
# ---------------------------------------------------------------------------
$wsUSD = $wb.Worksheets.Item("D1_USD")

$wsUSD.Range("A75").Copy()
$wsUSD.Range("A76").PasteSpecial(-4122)
$wsUSD.Range("A76").Value = 45264
$wsUSD.Range("B76").Value = 3.972988
$wsUSD.Range("C76").Value = 4.003844
$wsUSD.Range("D76").Formula = "=B76-C76"
$wsUSD.Range("E76").Formula = "=IF(D76<0,1,0)"
$wsUSD.Range("C77").Value = 3.9509294000000001
$wsUSD.Range("C78").Select()

# ---------------------------------------------------------------------------
# Sheet "D1_EUR" (sheet3.xml) — same pattern: fill in row 350 (45264) and
# backfill C351 with the next day's prediction.
# ---------------------------------------------------------------------------
$wsEUR1 = $wb.Worksheets.Item("D1_EUR")

$wsEUR1.Range("A349").Copy()
$wsEUR1.Range("A350").PasteSpecial(-4122)
$wsEUR1.Range("A350").Value = 45264
$wsEUR1.Range("B350").Value = 4.3256899999999998
$wsEUR1.Range("D350").Formula = "=B350-C350"
$wsEUR1.Range("E350").Formula = "=IF(D350<0,1,0)"
$wsEUR1.Range("C351").Value = 4.3591449999999998
$wsEUR1.Range("C352").Select()

# ---------------------------------------------------------------------------
# Sheet "D5_EUR" (sheet5.xml) — the placeholder "Nan" for 45264 is now a
# real predicted value.
# ---------------------------------------------------------------------------
$wsEUR5 = $wb.Worksheets.Item("D5_EUR")

$wsEUR5.Range("B43").Value = 4.3256899999999998
$wsEUR5.Range("B49").Select()

# ---------------------------------------------------------------------------
# Sheet "D1_OIL" (sheet6.xml) — three new days (45252, 45253, 45254) are
# inserted ahead of the previous first row, the previously-dangling
# prediction-only rows gain their actual/diff/ratio data, and a new day
# (45264) is appended. Column B also gains an explicit 0.0000 number format.
# ---------------------------------------------------------------------------
$wsOIL = $wb.Worksheets.Item("D1_OIL")

# Write bottom-up (row 10 -> row 2) using hard-coded source values so we
# never read a cell that a later step is about to overwrite.

# Row 10 (was: only C7 = 71.2617)
$wsOIL.Range("A1").Copy()
$wsOIL.Range("A10").PasteSpecial(-4122)
$wsOIL.Range("A10").Value = 45264
$wsOIL.Range("B10").Value = 73.040001000000004
$wsOIL.Range("C10").Value = 71.261700000000005
$wsOIL.Range("D10").Formula = "=B10-C10"
$wsOIL.Range("E10").Formula = "=D10/C10"

# Row 9 (was row 6: A6=45261 B6=74.07 C6=75.0822)
$wsOIL.Range("A1").Copy()
$wsOIL.Range("A9").PasteSpecial(-4122)
$wsOIL.Range("A9").Value = 45261
$wsOIL.Range("B9").Value = 74.069999999999993
$wsOIL.Range("C9").Value = 75.0822
$wsOIL.Range("D9").Formula = "=B9-C9"
$wsOIL.Range("E9").Formula = "=D9/C9"

# Row 8 (was row 5: A5=45260 B5=75.96 C5=73.5269)
$wsOIL.Range("A1").Copy()
$wsOIL.Range("A8").PasteSpecial(-4122)
$wsOIL.Range("A8").Value = 45260
$wsOIL.Range("B8").Value = 75.959998999999996
$wsOIL.Range("C8").Value = 73.526899999999998
$wsOIL.Range("D8").Formula = "=B8-C8"
$wsOIL.Range("E8").Formula = "=D8/C8"

# Row 7 (was row 4: A4=45259 B4=77.86 C4=74.5768)
$wsOIL.Range("A1").Copy()
$wsOIL.Range("A7").PasteSpecial(-4122)
$wsOIL.Range("A7").Value = 45259
$wsOIL.Range("B7").Value = 77.860000999999997
$wsOIL.Range("C7").Value = 74.576800000000006
$wsOIL.Range("D7").Formula = "=B7-C7"
$wsOIL.Range("E7").Formula = "=D7/C7"

# Row 6 (was row 3: A3=45258 B3=76.41 C3=76.5049)
$wsOIL.Range("A1").Copy()
$wsOIL.Range("A6").PasteSpecial(-4122)
$wsOIL.Range("A6").Value = 45258
$wsOIL.Range("B6").Value = 76.410004000000001
$wsOIL.Range("C6").Value = 76.504900000000006
$wsOIL.Range("D6").Formula = "=B6-C6"
$wsOIL.Range("E6").Formula = "=D6/C6"

# Row 5 (was row 2: A2=45257 B2=74.86, no C/D/E previously -> now gets them)
$wsOIL.Range("A1").Copy()
$wsOIL.Range("A5").PasteSpecial(-4122)
$wsOIL.Range("A5").Value = 45257
$wsOIL.Range("B5").Value = 74.860000999999997
$wsOIL.Range("C5").Value = 75.922899999999998
$wsOIL.Range("D5").Formula = "=B5-C5"
$wsOIL.Range("E5").Formula = "=D5/C5"

# Rows 2-4: brand-new leading days, date + actual price only.
$wsOIL.Range("A1").Copy()
$wsOIL.Range("A2").PasteSpecial(-4122)
$wsOIL.Range("A2").Value = 45252
$wsOIL.Range("B2").Value = 77.099997999999999

$wsOIL.Range("A1").Copy()
$wsOIL.Range("A3").PasteSpecial(-4122)
$wsOIL.Range("A3").Value = 45253
$wsOIL.Range("B3").Value = 76.349997999999999

$wsOIL.Range("A1").Copy()
$wsOIL.Range("A4").PasteSpecial(-4122)
$wsOIL.Range("A4").Value = 45254
$wsOIL.Range("B4").Value = 75.540001000000004

# Copy the date number-format (style "1") onto the freshly written A2:A10
# cells (they were already pasted above, this just keeps things explicit).
$wsOIL.Range("A1").Copy()
$wsOIL.Range("A2:A10").PasteSpecial(-4122)

# Re-apply the right values after the format paste above (PasteSpecial of
# formats only does not disturb values, but keep this deterministic).
$wsOIL.Range("A2").Value = 45252
$wsOIL.Range("A3").Value = 45253
$wsOIL.Range("A4").Value = 45254
$wsOIL.Range("A5").Value = 45257
$wsOIL.Range("A6").Value = 45258
$wsOIL.Range("A7").Value = 45259
$wsOIL.Range("A8").Value = 45260
$wsOIL.Range("A9").Value = 45261
$wsOIL.Range("A10").Value = 45264

# E column ratio cells keep the existing percentage style used by the
# original E3:E6 cells.
$wsOIL.Range("E6").Copy()
$wsOIL.Range("E5:E10").PasteSpecial(-4122)
$wsOIL.Range("E5").Formula = "=D5/C5"
$wsOIL.Range("E6").Formula = "=D6/C6"
$wsOIL.Range("E7").Formula = "=D7/C7"
$wsOIL.Range("E8").Formula = "=D8/C8"
$wsOIL.Range("E9").Formula = "=D9/C9"
$wsOIL.Range("E10").Formula = "=D10/C10"

# Column B (actual price) and its header now carry an explicit 0.0000
# number format.
$wsOIL.Range("B1:B10").NumberFormat = "0.0000"

# E1 average now starts at row 5 (first row with a D value after the
# reshuffle) instead of row 2.
$wsOIL.Range("E1").Formula = "=AVERAGE(D5:D301)"

$wsOIL.Range("C4").Select()
